$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Women")

$ws.Range("A8").Value = 16
$ws.Range("B8").Value = "uhoujä"
$ws.Range("C8").Value = "höjhk"
